$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the drop-table shared strings so probability values are scaled
# from fractional (0-1) to a 10000-based integer range (batch/stage import format).
$ws.Range("G3").Value = "(type:Experience;amount:10;probability:10000;guaranteed:true)(type:Gold;amount:5;probability:10000;guaranteed:true)(type:Health;amount:5;probability:1000;guaranteed:false)(type:Ability;amount:1;probability:500;guaranteed:false)"
$ws.Range("G4").Value = "(type:Ability;amount:1;probability:10000;guaranteed:false)"
$ws.Range("G5").Value = "(type:Experience;amount:10;probability:10000;guaranteed:true)(type:Gold;amount:5;probability:10000;guaranteed:true)(type:Health;amount:5;probability:1000;guaranteed:false)(type:Ability;amount:1;probability:500;guaranteed:false)"

# Scale up the numeric stat columns (C:F) for rows 3-5 by 10000.
$ws.Range("C3").Value = 500000
$ws.Range("D3").Value = 50000
$ws.Range("E3").Value = 15000
$ws.Range("F3").Value = 5000

$ws.Range("C4").Value = 500000
$ws.Range("D4").Value = 50000
$ws.Range("E4").Value = 15000
$ws.Range("F4").Value = 10000

$ws.Range("C5").Value = 500000
$ws.Range("D5").Value = 50000
$ws.Range("E5").Value = 15000
$ws.Range("F5").Value = 10000

# Move the active selection to G5, matching the saved view state.
$ws.Range("G5").Select()
